$wb = $excel.ActiveWorkbook

# --- Update comments on "Sample Section" header row (A1:V1) ---
$wsMain = $wb.Worksheets.Item("Sample Section")

[void]$wsMain.Range("A1").Comment.Text("(Required) The unique identifier from HuBMAP or SenNet for the source (parent`ndata) from which the sample was derived. Example: HBM122.EFGH.789")
[void]$wsMain.Range("B1").Comment.Text("(Required) The unique HuBMAP or SenNet identifier assigned to the sample by the`ningest portal. Example: HBM743.CKJW.876")
[void]$wsMain.Range("C1").Comment.Text("A locally assigned identifier provided by the data provider for the dataset. It`nis used to reference an external metadata record that may be maintained`nindependently, enabling traceability and supporting provenance tracking.`nExample: Visium_9OLC_A4_S1")
[void]$wsMain.Range("D1").Comment.Text("(Required) The DOI for the protocols.io page that details the assay or the`nprocedures used for sample procurement and preparation. For example, in the case`nof an imaging assay, the protocol may start with tissue section staining and end`nwith the generation of an OME-TIFF file. The documented protocol should also`ninclude any image processing steps involved in producing the final OME-TIFF.`nExample: https://dx.doi.org/10.17504/protocols.io.eq2lyno9qvx9/v1")
[void]$wsMain.Range("E1").Comment.Text("(Required) The length of time the sample was stored prior to processing it. For`nassays performed on tissue sections, this refers to how long the tissue section`n(e.g., slide) was stored before the assay began (e.g., imaging). For assays`nperformed on suspensions, such as sequencing, it refers to how long the`nsuspension was stored before library construction started. Example: 12")
[void]$wsMain.Range("F1").Comment.Text("(Required) The unit of measurement used to specify the source storage duration`nvalue. Example: hour")
[void]$wsMain.Range("G1").Comment.Text("(Required) The medium used during the sample preparation process. If no specific`nmedium was utilized, enter `"None`". If medium was not recorded, enter `"Unknown`".`nExample: Fresh frozen CMC")
[void]$wsMain.Range("H1").Comment.Text("(Required) The condition under which the sample preparation took place, such as`nwhether the sample was placed on dry ice during the process. If preparation`ncondition was not recorded, enter `"Unknown`". Example: Frozen on dry ice")
[void]$wsMain.Range("I1").Comment.Text("The duration for which the tissue was handled prior to its initial preservation.`nExample: 120")
[void]$wsMain.Range("J1").Comment.Text("The unit of measurement for the processing time value. If processing time is not`nspecified, this field may be left blank. Example: minute")
[void]$wsMain.Range("K1").Comment.Text("(Required) The medium used to preserve the sample. If no specific medium was`nutilized, enter `"None`". If medium was not recorded, enter `"Unknown`". Example:`nFFPE (Paraffin embedded)")
[void]$wsMain.Range("L1").Comment.Text("(Required) The method used to store the sample after preparation and prior to`nperforming the assay. If no specific storage method was utilized, enter `"None`".`nIf storage method was not recorded, enter `"Unknown`". Example: Frozen in dry ice")
[void]$wsMain.Range("M1").Comment.Text("The quality criteria used to assess the sample, which may include metrics such`nas RIN (e.g., RIN: 8.7) or visual inspection parameters for suspensions prior to`ncell lysis. These criteria can be captured at a high level with general terms`nlike `"OK`" or `"not OK`" or with more specific descriptors such as `"debris`" `"clump`"`nor `"low clump`". Example: RIN: 8.7, low clump, no visible debris")
[void]$wsMain.Range("N1").Comment.Text("The key variables in the histopathological report that are crucial for assessing`nthe tissue, including the absence of necrosis, comments on tissue composition,`ndescriptions of significant pathology, and high-level assessments of`ninflammation or fibrosis. Example: No necrosis observed; tissue composed`npredominantly of hepatocytes with mild portal inflammation and minimal fibrosis")
[void]$wsMain.Range("O1").Comment.Text("(Required) The thickness of an object in question. Example: 10")
[void]$wsMain.Range("P1").Comment.Text("(Required) The unit of measurement for the thickness value. If no thickness`nmeasurement is specified, this field may be left blank. Example: mm")
[void]$wsMain.Range("Q1").Comment.Text("(Required) The index number assigned to the tissue section, with numbering`nbeginning at 1 for sections within a block. Example: 1")
[void]$wsMain.Range("R1").Comment.Text("The area of the object being measured. Example: 100")
[void]$wsMain.Range("S1").Comment.Text("The unit of measurement used to define the area. If no area value is specified,`nthis field may be left blank. Example: mm^2")
[void]$wsMain.Range("T1").Comment.Text("Indicates whether the section was rehydrated. Example: No")
[void]$wsMain.Range("U1").Comment.Text("Miscellaneous details about the sample that are not captured in the existing`nmetadata fields. Example: Sample was stored at 4°C for 48 hours prior to`nprocessing due to equipment maintenance delay")
[void]$wsMain.Range("V1").Comment.Text("(Required) The unique string identifier for the metadata specification version,`nwhich is easily interpretable by computers for purposes of data validation and`nprocessing. Example: 22bc762a-5020-419d-b170-24253ed9e8d9")

# --- Update "storage_medium" lookup sheet: add "Water" option and reorder rows ---
$wsStorage = $wb.Worksheets.Item("storage_medium")

$wsStorage.Cells.Item(1, 1).Value = "Water"
$wsStorage.Cells.Item(1, 2).Value = "http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C65147"
$wsStorage.Cells.Item(2, 1).Value = "OCT"
$wsStorage.Cells.Item(2, 2).Value = "http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C63523"
$wsStorage.Cells.Item(3, 1).Value = "NBF (Neutral Buffered Formalin)"
$wsStorage.Cells.Item(3, 2).Value = "http://purl.obolibrary.org/obo/OBIB_0000213"
$wsStorage.Cells.Item(4, 1).Value = "Allprotect tissue reagent (ALL)"
$wsStorage.Cells.Item(4, 2).Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000118"
$wsStorage.Cells.Item(5, 1).Value = "DMSO (no serum)"
$wsStorage.Cells.Item(5, 2).Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000115"
$wsStorage.Cells.Item(6, 1).Value = "PFA (Paraformaldehyde)"
$wsStorage.Cells.Item(6, 2).Value = "http://purl.obolibrary.org/obo/CHEBI_61538"
$wsStorage.Cells.Item(7, 1).Value = "Unknown"
$wsStorage.Cells.Item(7, 2).Value = "http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C17998"
$wsStorage.Cells.Item(8, 1).Value = "Gelatin"
$wsStorage.Cells.Item(8, 2).Value = "http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C65802"
$wsStorage.Cells.Item(9, 1).Value = "DMSO (serum)"
$wsStorage.Cells.Item(9, 2).Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000125"
$wsStorage.Cells.Item(10, 1).Value = "CMC"
$wsStorage.Cells.Item(10, 2).Value = "http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C83594"
$wsStorage.Cells.Item(11, 1).Value = "2% PFA/2.5% Glutaraldehyde"
$wsStorage.Cells.Item(11, 2).Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000332"
$wsStorage.Cells.Item(12, 1).Value = "Methanol"
$wsStorage.Cells.Item(12, 2).Value = "http://purl.obolibrary.org/obo/CHEBI_17790"
$wsStorage.Cells.Item(13, 1).Value = "PAXgene tissue kit (PXT)"
$wsStorage.Cells.Item(13, 2).Value = "http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C185113"
$wsStorage.Cells.Item(14, 1).Value = "PBS"
$wsStorage.Cells.Item(14, 2).Value = "http://purl.obolibrary.org/obo/OBI_0100046"
$wsStorage.Cells.Item(15, 1).Value = "1X quench buffer"
$wsStorage.Cells.Item(15, 2).Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000427"
$wsStorage.Cells.Item(16, 1).Value = "Ethanol"
$wsStorage.Cells.Item(16, 2).Value = "http://purl.obolibrary.org/obo/CHEBI_16236"
$wsStorage.Cells.Item(17, 1).Value = "Formic acid in water"
$wsStorage.Cells.Item(17, 2).Value = "http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C83719"
$wsStorage.Cells.Item(18, 1).Value = "HPMC-PVP"
$wsStorage.Cells.Item(18, 2).Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000386"
$wsStorage.Cells.Item(19, 1).Value = "MACS tissue storage solution"
$wsStorage.Cells.Item(19, 2).Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000105"
$wsStorage.Cells.Item(20, 1).Value = "Tris-EDTA"
$wsStorage.Cells.Item(20, 2).Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000135"
$wsStorage.Cells.Item(21, 1).Value = "Concentrated quench buffer"
$wsStorage.Cells.Item(21, 2).Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000391"
$wsStorage.Cells.Item(22, 1).Value = "Cryo-EM"
$wsStorage.Cells.Item(22, 2).Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000333"
$wsStorage.Cells.Item(23, 1).Value = "RNAlater"
$wsStorage.Cells.Item(23, 2).Value = "http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C63348"
$wsStorage.Cells.Item(24, 1).Value = "FFPE (Paraffin embedded)"
$wsStorage.Cells.Item(24, 2).Value = "http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C143028"
$wsStorage.Cells.Item(25, 1).Value = "None"
$wsStorage.Cells.Item(25, 2).Value = "http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C41132"

# --- Expand the K column (storage_medium) data validation range to include new row ---
$wsMain.Range("K2:K1001").Validation.Formula1 = "'storage_medium'!`$A`$1:`$A`$25"

# --- Update pav:createdOn timestamp on the .metadata sheet ---
$wsMeta = $wb.Worksheets.Item(".metadata")
$wsMeta.Range("C2").Value = "2025-10-16T07:27:11-07:00"
